$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.98"
$ws.Range("E2").Value = "'0.38%"
$ws.Range("D3").Value = "'41.08"
$ws.Range("E3").Value = "'-0.21%"
$ws.Range("D4").Value = "'5.209"
$ws.Range("E4").Value = "'1.74%"
$ws.Range("D5").Value = "'0.07663"
$ws.Range("E5").Value = "'0.26%"
$ws.Range("D6").Value = "'1.747"
$ws.Range("E6").Value = "'7.78%"
$ws.Range("D7").Value = "'0.9210"
$ws.Range("E7").Value = "'1.65%"
$ws.Range("E8").Value = "'-1.22%"
$ws.Range("D9").Value = "'0.1269"
$ws.Range("E9").Value = "'12.77%"
$ws.Range("D10").Value = "'0.1821"
$ws.Range("E10").Value = "'0.90%"
$ws.Range("D11").Value = "'0.09099"
$ws.Range("E11").Value = "'-0.14%"
$ws.Range("D12").Value = "'0.04159"
$ws.Range("E12").Value = "'-2.07%"
$ws.Range("E13").Value = "'-0.01%"
$ws.Range("D14").Value = "'0.001284"
$ws.Range("E14").Value = "'2.78%"
$ws.Range("D15").Value = "'0.005895"
$ws.Range("E15").Value = "'1.68%"
$ws.Range("D16").Value = "'3.353"
$ws.Range("E16").Value = "'0.21%"
$ws.Range("D17").Value = "'4.298"
$ws.Range("E17").Value = "'0.66%"
$ws.Range("D19").Value = "'7.389"
$ws.Range("E19").Value = "'9.67%"
$ws.Range("D20").Value = "'0.1354"
$ws.Range("E20").Value = "'-0.41%"
$ws.Range("D21").Value = "'0.2723"
$ws.Range("E21").Value = "'-0.45%"
$ws.Range("D22").Value = "'0.04020"
$ws.Range("E22").Value = "'-1.08%"
$ws.Range("E23").Value = "'0.24%"
$ws.Range("D24").Value = "'0.004095"
$ws.Range("E24").Value = "'1.37%"
$ws.Range("D25").Value = "'0.0001272"
$ws.Range("E25").Value = "'0.07%"
$ws.Range("D38").Value = "'0.02511"
$ws.Range("E38").Value = "'3.59%"
$ws.Range("D39").Value = "'0.05313"
$ws.Range("E39").Value = "'1.22%"
$ws.Range("D40").Value = "'0.007856"
$ws.Range("E40").Value = "'0.59%"
$ws.Range("E41").Value = "'0.71%"
$ws.Range("D42").Value = "'0.006656"
$ws.Range("E42").Value = "'1.96%"
$ws.Range("D43").Value = "'0.002054"
$ws.Range("E43").Value = "'5.26%"
$ws.Range("D44").Value = "'0.008122"
$ws.Range("E44").Value = "'7.13%"
$ws.Range("D45").Value = "'0.3081"
$ws.Range("E45").Value = "'-0.11%"
$ws.Range("D46").Value = "'0.00006798"
$ws.Range("E46").Value = "'0.31%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("D48").Value = "'0.2241"
$ws.Range("E48").Value = "'256.50%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.14%"
